$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 2.22
$ws.Range("L2").Value = 1.44
$ws.Range("R2").Value = 1.33
$ws.Range("T2").Value = 1.83
$ws.Range("V2").Value = 1.81
$ws.Range("Y2").Value = 9.199999999999999
$ws.Range("Z2").Value = 16
$ws.Range("AB2").Value = 15
$ws.Range("AD2").Value = 12.5
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 48
$ws.Range("AK2").Value = 55
$ws.Range("AL2").Value = 60
$ws.Range("AN2").Value = 65

# Row 4 updates
$ws.Range("N4").Value = 1.32
$ws.Range("P4").Value = 1.32

# Row 7 updates
$ws.Range("G7").Value = 3.75
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 2.56
$ws.Range("J7").Value = 3.8
$ws.Range("P7").Value = 2.44
$ws.Range("Q7").Value = 1.43

# Row 8 updates
$ws.Range("H8").Value = 2.82
$ws.Range("K8").Value = 7.6
$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 1.48

# Row 9 updates
$ws.Range("I9").Value = 1.96
$ws.Range("Q9").Value = 1.93

# Row 10 updates
$ws.Range("P10").Value = 3.15

# Row 12 updates
$ws.Range("Q12").Value = 1.41
